# Auto-update draw results: append the 2025-10-25 Pick 4 result as a new
# row at the bottom of the "Results" sheet (mirrors the nightly scraper
# commit that appends one row per draw).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the current last populated row in column A, then target the row
# right after it for the new record.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$targetRow = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 5))
$lastRowRange = $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, 5))

# Force text storage first so date-like / numeric-like strings (e.g. the
# draw date and the 6-digit phase code) are kept as literal text instead
# of being parsed into a date serial number or a numeric value, matching
# how every other row in the sheet is stored.
$targetRow.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2025-10-25"
$ws.Cells.Item($newRow, 2).Value = "Pick 4"
$ws.Cells.Item($newRow, 3).Value = "251025"
$ws.Cells.Item($newRow, 4).Value = "7-9-6-3"
$ws.Cells.Item($newRow, 5).Value = "2025-10-25T21:36:06.852+04:00"

# Re-apply the same style as the preceding row so the new row doesn't
# pick up a distinct "Text" cell style (keeps formatting consistent with
# the rest of the table).
$targetRow.Style = $lastRowRange.Style
